$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A26").NumberFormat = "@"
$ws.Range("A26").Value = "02/06/2026"
$ws.Range("A26").Style = "Normal"

$ws.Range("B26").Value = 1767.745999999999
$ws.Range("C26").Value = 0.0280017604339085
$ws.Range("D26").Value = 50
